$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("Abstenção") - this shifts the
# existing Gravidade..Observação columns from B:M to C:N, carrying the
# header style (bold/border) along for the new column and for every
# shifted header cell.
$ws.Columns.Item(2).Insert()

# New header for the inserted column.
$ws.Range("B1").Value = "Abstenção"

# Row 2 now holds a new evaluator submission - overwrite every data cell
# with the final values (the old M2 "Observação" note, shifted to N2 by
# the column insert, is no longer present so it is cleared).
$ws.Range("A2").Value = "Usuário Teste"
$ws.Range("B2").Value = "Não"
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 4.5
$ws.Range("E2").Value = 4.5
$ws.Range("F2").Value = 3.75
$ws.Range("G2").Value = 3.5
$ws.Range("H2").Value = 4
$ws.Range("I2").Value = 4.5
$ws.Range("J2").Value = 3.5
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 3.8
$ws.Range("M2").Value = 3.775
$ws.Range("N2").ClearContents()
